$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.493.86"
$ws.Range("D3").Value = "2.067.11"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.76"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.37%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0784"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "2.371.12"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.763"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "2.073.82"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "37.421.79"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("D21").Value = "0.0$([char]0x2083)0831"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("E28").Value = "  -5.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.51%  "
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("E31").Value = "  -0.51%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.47%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.31"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0227"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.28%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.88%  "
$ws.Range("D43").Value = "1.491.61"
$ws.Range("E43").Value = "  +2.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0954"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.86%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.04"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.03"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.27"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("D51").Value = "2.256.07"
$ws.Range("E51").Value = "  -0.95%  "
